$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header column "fecha"
$ws.Cells.Item(1, 8).Value = "fecha"

# Data rows: A, B, C, D, E, F, G, H(fecha - kept as text)
$data = @(
    @(1,  1, 234, 432, 65, 65, 123000, "2020-07-12"),
    @(6,  1, 492, 485,  7,  7, 165200, "2020-01-25"),
    @(7,  1, 492, 485,  7,  7,  43500, "2020-01-13"),
    @(8,  1, 492, 485,  7,  7,  43500, "2020-01-13"),
    @(9,  1, 492, 485,  7,  7,  43500, "2020-01-13"),
    @(10, 1, 492, 485,  7,  7,  43500, "2020-01-13"),
    @(13, 1, 492, 485,  7,  7,  43500, "2020-01-13"),
    @(15, 1, 492, 485,  7,  7,  43500, "2020-08-30"),
    @(16, 1, 492, 485,  7,  7,  43500, "2020-08-30")
)

# Keep the fecha column as text so dates like "2020-07-12" aren't
# reinterpreted as date serial numbers.
$ws.Range("H2:H10").NumberFormat = "@"

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $ws.Cells.Item($rowIndex, 6).Value = $row[5]
    $ws.Cells.Item($rowIndex, 7).Value = $row[6]
    $ws.Cells.Item($rowIndex, 8).Value = $row[7]
    $rowIndex++
}
